# "Add files via upload" - refresh the indicator metadata sheet (13-1-1)
# with the Kyrgyz Republic National Statistical Committee's current
# contact / data-reporter details.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "2. Data reporter" block (rows 6-10).
# Set in the same order the new shared strings are appended in the target
# workbook (website, phone, contact person, organization) so the rebuilt
# sharedStrings table lines up.
$ws.Range("B10").Value = "www.stat.gov.kg"
$ws.Range("B9").Value  = "(0312) 62 56 07"
$ws.Range("B7").Value  = "Mambetaliev T.A."
$ws.Range("B6").Value  = "The National Statistical Committee of the Kyrgyz Republic (Department of Digital Development and Sustainable Development Statistics)"

# Move the active selection to the contact-person cell.
$ws.Range("B7").Select() | Out-Null
